# Generate Report for Handoff
# Refresh the "Latest Handoff"/"Latest Handback" timestamps for the most
# recently processed file (6691b810-1c0d-4806-b331-2af060662bd0) across the
# Overview summary sheet and each per-locale detail sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D7").Value = "2016-28-12 10:28:47"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E7").Value = "2016-03-12 10:28:44"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E7").Value = "2016-03-12 10:28:47"
